$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.9008665625316319, 0.6854350136727263, 0.6841829112431119, 3.2037666687489232)"
$ws.Range("C2").Value = "NIG(1.7811383103389349, 1.35683092714103, 2.827121312230168, 6.860218842387461)"
$ws.Range("D2").Value = "NIG(0.9003714525874157, 0.628067454188995, 1.97043778914033, 2.9684966470693817)"
$ws.Range("E2").Value = "NIG(2.0656295765146866, 1.7350338020771925, 2.8783357149566995, 5.671768437567078)"
